# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# D.Prescott (row 2)
$rushing.Cells.Item(2, 3).Value = 12   # 1DATT
$rushing.Cells.Item(2, 4).Value = 5    # 2DATT
$rushing.Cells.Item(2, 5).Value = 11   # 3DATT

# C.Rush (row 3)
$rushing.Cells.Item(3, 3).Value = 2    # 1DATT

# E.Elliott (row 4)
$rushing.Cells.Item(4, 3).Value = 125  # 1DATT
$rushing.Cells.Item(4, 4).Value = 59   # 2DATT
$rushing.Cells.Item(4, 5).Value = 27   # 3DATT
$rushing.Cells.Item(4, 6).Value = 38   # RZATT

# T.Pollard (row 5)
$rushing.Cells.Item(5, 3).Value = 76   # 1DATT
$rushing.Cells.Item(5, 4).Value = 42   # 2DATT
$rushing.Cells.Item(5, 5).Value = 8    # 3DATT
$rushing.Cells.Item(5, 6).Value = 16   # RZATT

# C.Clement (row 6)
$rushing.Cells.Item(6, 3).Value = 11   # 1DATT
$rushing.Cells.Item(6, 4).Value = 5    # 2DATT
$rushing.Cells.Item(6, 5).Value = 4    # 3DATT

# C.Lamb (row 8)
$rushing.Cells.Item(8, 3).Value = 4    # 1DATT

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# E.Elliott (row 2)
$receiving.Cells.Item(2, 3).Value = 45   # Short Target
$receiving.Cells.Item(2, 4).Value = 33   # Short Comp
$receiving.Cells.Item(2, 7).Value = 17   # RZ Target
$receiving.Cells.Item(2, 8).Value = 12   # RZ Comp

# T.Pollard (row 3)
$receiving.Cells.Item(3, 3).Value = 32   # Short Target
$receiving.Cells.Item(3, 4).Value = 27   # Short Comp
$receiving.Cells.Item(3, 7).Value = 4    # RZ Target

# A.Cooper (row 5)
$receiving.Cells.Item(5, 3).Value = 65   # Short Target
$receiving.Cells.Item(5, 4).Value = 48   # Short Comp
$receiving.Cells.Item(5, 5).Value = 25   # Deep Target
$receiving.Cells.Item(5, 6).Value = 13   # Deep Comp
$receiving.Cells.Item(5, 7).Value = 17   # RZ Target
$receiving.Cells.Item(5, 8).Value = 13   # RZ Comp

# C.Lamb (row 6)
$receiving.Cells.Item(6, 3).Value = 76   # Short Target
$receiving.Cells.Item(6, 4).Value = 53   # Short Comp
$receiving.Cells.Item(6, 5).Value = 34   # Deep Target
$receiving.Cells.Item(6, 6).Value = 18   # Deep Comp

# M.Gallup (row 7)
$receiving.Cells.Item(7, 3).Value = 44   # Short Target
$receiving.Cells.Item(7, 4).Value = 34   # Short Comp
$receiving.Cells.Item(7, 5).Value = 14   # Deep Target
$receiving.Cells.Item(7, 6).Value = 8    # Deep Comp
$receiving.Cells.Item(7, 7).Value = 6    # RZ Target
$receiving.Cells.Item(7, 8).Value = 4    # RZ Comp

# C.Wilson (row 8)
$receiving.Cells.Item(8, 3).Value = 24   # Short Target
$receiving.Cells.Item(8, 4).Value = 19   # Short Comp
$receiving.Cells.Item(8, 5).Value = 10   # Deep Target

# M.Turner (row 10)
$receiving.Cells.Item(10, 3).Value = 11  # Short Target
$receiving.Cells.Item(10, 4).Value = 10  # Short Comp
$receiving.Cells.Item(10, 5).Value = 3   # Deep Target
$receiving.Cells.Item(10, 6).Value = 1   # Deep Comp
$receiving.Cells.Item(10, 7).Value = 5   # RZ Target
$receiving.Cells.Item(10, 8).Value = 4   # RZ Comp

# D.Schultz (row 12)
$receiving.Cells.Item(12, 3).Value = 67  # Short Target
$receiving.Cells.Item(12, 4).Value = 55  # Short Comp
$receiving.Cells.Item(12, 5).Value = 9   # Deep Target
$receiving.Cells.Item(12, 7).Value = 11  # RZ Target
$receiving.Cells.Item(12, 8).Value = 7   # RZ Comp
